# Revert "Updated documentation and climates"
# This undoes the VPDF-row formatting/units that had been added to the
# "325_-35" sheet (row 5: B5/C5/D5/E5/F5), restoring it to match the
# still-unedited "275_-45" sheet's row 5.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "275_-45" - untouched reference sheet
$ws2 = $wb.Worksheets.Item(2)   # "325_-35" - sheet that was edited

# B5 had been given a new "mm" label with a new font/style (s=21); put the
# original formatting back (copied from the matching, never-edited cell on
# sheet 1) and clear its value.
$ws1.Range("B5").Copy()
$ws2.Range("B5").PasteSpecial(-4122)
$ws2.Range("B5").Value = ""

# C5:E5 had their unit label changed from "cm.cm-1" to "mm.mm-1" - restore it.
$ws2.Range("C5:E5").Value = "cm.cm-1"

# F5 had picked up a stray "mm.mm-1" label - it should be blank again.
$ws2.Range("F5").Value = ""

# Restore the previously-active selection on the "325_-35" sheet.
$ws2.Activate() | Out-Null
$ws2.Range("H2").Select() | Out-Null
